$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.087.73'
$ws.Range('E2').Value = '  +0.40%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.325.56'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '97.95'
$ws.Range('E5').Value = '  +3.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '272.05'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.631'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.628'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.58'
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.02'
$ws.Range('E12').Value = '  -3.77%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.662.59'
$ws.Range('E14').Value = '  +3.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.54'
$ws.Range('E15').Value = '  +1.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.875'
$ws.Range('E16').Value = '  +7.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.322.33'
$ws.Range('E17').Value = '  +3.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.999.62'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('E19').Value = '  +4.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.42'
$ws.Range('E20').Value = '  +4.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.59'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.44'
$ws.Range('E23').Value = '  +2.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.40'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  +1.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.43'
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  -2.31%  '
$ws.Range('E29').Value = '  +2.07%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.36'
$ws.Range('E30').Value = '  -5.16%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.43'
$ws.Range('E31').Value = '  +6.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '175.59'
$ws.Range('E32').Value = '  +1.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0915'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.52'
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  +2.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0366'
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.47'
$ws.Range('E38').Value = '  +3.32%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.39'
$ws.Range('E39').Value = '  -5.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.245'
$ws.Range('E40').Value = '  +7.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.39'
$ws.Range('E41').Value = '  +10.28%  '
$ws.Range('E42').Value = '  +22.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.43'
$ws.Range('E43').Value = '  -3.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '62.99'
$ws.Range('E44').Value = '  -0.90%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '9.12'
$ws.Range('E45').Value = '  +8.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.35'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('E47').Value = '  +3.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '100.52'
$ws.Range('E48').Value = '  -0.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.21'
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +16.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.549.43'
$ws.Range('E51').Value = '  +3.55%  '
